$wb = $excel.ActiveWorkbook

# ===== Section_A =====
$ws = $wb.Worksheets.Item("Section_A")
$ws.Range("B2").Value = "MA262 [C304]"
$ws.Range("C2").Value = "CS261 [C002]"
$ws.Range("D2").Value = "Free"
$ws.Range("E2").Value = "CS262 [C403]"
$ws.Range("C3").Value = "CS263 [C303]"
$ws.Range("D3").Value = "CS264 [C201]"
$ws.Range("E3").Value = "Free"
$ws.Range("F3").Value = "CS264 [C201]"
$ws.Range("B5").Value = "ELECTIVE_B3 [C204]"
$ws.Range("C5").Value = "MA261 [C401]"
$ws.Range("D5").Value = "ELECTIVE_B3 [C204]"
$ws.Range("E5").Value = "CS263 [C303]"
$ws.Range("F5").Value = "CS262 [C403]"
$ws.Range("B6").Value = "Free"
$ws.Range("C6").Value = "Free"
$ws.Range("B7").Value = "MA261 [C401]"
$ws.Range("C7").Value = "CS263 (Lab) [L306]"
$ws.Range("D7").Value = "CS262 (Lab) [L408]"
$ws.Range("E7").Value = "MA262 [C304]"
$ws.Range("F7").Value = "CS261 [C002]"
$ws.Range("B8").Value = "CS264 (Tutorial) [C102]"
$ws.Range("C8").Value = "CS263 (Lab) [L306]"
$ws.Range("D8").Value = "CS262 (Lab) [L408]"
$ws.Range("E8").Value = "Free"

# ===== Section_B =====
$ws = $wb.Worksheets.Item("Section_B")
$ws.Range("B2").Value = "Free"
$ws.Range("C2").Value = "MA261 [C201]"
$ws.Range("D2").Value = "MA261 [C201]"
$ws.Range("E2").Value = "Free"
$ws.Range("F2").Value = "CS263 [C102]"
$ws.Range("B3").Value = "CS261 [C003]"
$ws.Range("C3").Value = "CS263 [C102]"
$ws.Range("D3").Value = "MA262 [C202]"
$ws.Range("F3").Value = "CS264 [C305]"
$ws.Range("B5").Value = "ELECTIVE_B3 [C104]"
$ws.Range("C5").Value = "CS264 [C305]"
$ws.Range("D5").Value = "ELECTIVE_B3 [C104]"
$ws.Range("E5").Value = "MA262 [C202]"
$ws.Range("F5").Value = "CS261 [C003]"
$ws.Range("D6").Value = "CS264 (Tutorial) [C204]"
$ws.Range("F6").Value = "Free"
$ws.Range("B7").Value = "CS262 [C201]"
$ws.Range("C7").Value = "Free"
$ws.Range("D7").Value = "CS263 (Lab) [L406]"
$ws.Range("E7").Value = "CS262 [C201]"
$ws.Range("F7").Value = "CS262 (Lab) [L408]"
$ws.Range("C8").Value = "Free"
$ws.Range("D8").Value = "CS263 (Lab) [L406]"
$ws.Range("E8").Value = "Free"
$ws.Range("F8").Value = "CS262 (Lab) [L408]"

# ===== Verification_A =====
$ws = $wb.Worksheets.Item("Verification_A")
$ws.Range("A2").Value = "**MA262**"
$ws.Range("B2").Value = "Differential Equations"
$ws.Range("C2").Value = "Anand Barangi"
$ws.Range("D2").Value = "3-1-0-0-2"
$ws.Range("F2").Value = "0/0"
$ws.Range("I2").Value = "C304"
$ws.Range("I3").Value = "C204"
$ws.Range("A4").Value = "**MA261**"
$ws.Range("B4").Value = "Multivariable Calculus"
$ws.Range("C4").Value = "Somen B"
$ws.Range("D4").Value = "3-1-0-0-2"
$ws.Range("E4").Value = "2/0"
$ws.Range("I4").Value = "C401"
$ws.Range("A5").Value = "**CS264**"
$ws.Range("B5").Value = "Computer Networks"
$ws.Range("C5").Value = "Prabhu Prasad B M"
$ws.Range("D5").Value = "3-1-0-0-4"
$ws.Range("E5").Value = "2/1"
$ws.Range("I5").Value = "C102, C201"
$ws.Range("A6").Value = "**CS261**"
$ws.Range("B6").Value = "Operating systems"
$ws.Range("C6").Value = "Suvadip Hazra"
$ws.Range("D6").Value = "3-0-0-4-2"
$ws.Range("F6").Value = "0/0"
$ws.Range("H6").Value = "Partial"
$ws.Range("I6").Value = "C002"
$ws.Range("A7").Value = "**CS263**"
$ws.Range("B7").Value = "Design & Analysis of Algorithms"
$ws.Range("C7").Value = "Malay, Pramod Y"
$ws.Range("D7").Value = "3-0-2-0-4"
$ws.Range("F7").Value = "2/1"
$ws.Range("I7").Value = "C303, L306"
$ws.Range("A8").Value = "**CS262**"
$ws.Range("B8").Value = "Software design tool and tecnique"
$ws.Range("C8").Value = "Sunil P V, Vivekraj"
$ws.Range("D8").Value = "2-0-2-0-3"
$ws.Range("F8").Value = "2/1"
$ws.Range("H8").Value = "Complete"
$ws.Range("I8").Value = "C403, L408"
$ws.Range("H9").Value = "[WARN] 6 issues"

# ===== Verification_B =====
$ws = $wb.Worksheets.Item("Verification_B")
$ws.Range("A2").Value = "**CS261**"
$ws.Range("B2").Value = "Operating systems"
$ws.Range("C2").Value = "Suvadip Hazra"
$ws.Range("D2").Value = "3-0-0-4-2"
$ws.Range("E2").Value = "2/0"
$ws.Range("I2").Value = "C003"
$ws.Range("A3").Value = "**ELECTIVE_B3**"
$ws.Range("B3").Value = "Elective Basket"
$ws.Range("C3").Value = "–"
$ws.Range("D3").Value = "3-0-0-0-3"
$ws.Range("E3").Value = "0/0"
$ws.Range("I3").Value = "C104"
$ws.Range("A4").Value = "**CS262**"
$ws.Range("B4").Value = "Software design tool and tecnique"
$ws.Range("C4").Value = "Sunil P V, Vivekraj"
$ws.Range("D4").Value = "2-0-2-0-3"
$ws.Range("E4").Value = "2/0"
$ws.Range("F4").Value = "2/1"
$ws.Range("H4").Value = "Complete"
$ws.Range("I4").Value = "L408, C201"
$ws.Range("A5").Value = "**MA261**"
$ws.Range("B5").Value = "Multivariable Calculus"
$ws.Range("C5").Value = "Somen B"
$ws.Range("D5").Value = "3-1-0-0-2"
$ws.Range("F5").Value = "0/0"
$ws.Range("I5").Value = "C201"
$ws.Range("A6").Value = "**CS263**"
$ws.Range("B6").Value = "Design & Analysis of Algorithms"
$ws.Range("C6").Value = "Malay, Pramod Y"
$ws.Range("D6").Value = "3-0-2-0-4"
$ws.Range("H6").Value = "Partial"
$ws.Range("I6").Value = "C102, L406"
$ws.Range("A7").Value = "**CS264**"
$ws.Range("B7").Value = "Computer Networks"
$ws.Range("C7").Value = "Prabhu Prasad B M"
$ws.Range("D7").Value = "3-1-0-0-4"
$ws.Range("E7").Value = "2/1"
$ws.Range("I7").Value = "C305, C204"
$ws.Range("I8").Value = "C202"
$ws.Range("H9").Value = "[WARN] 6 issues"

# ===== Room_Allocation =====
$ws = $wb.Worksheets.Item("Room_Allocation")
$ws.Range("A1").Value = "Room Number"
$ws.Range("B1").Value = "Type"
$ws.Range("C1").Value = "Capacity"
$ws.Range("D1").Value = "Facilities"
$ws.Range("E1").Value = "Total Sessions"
$ws.Range("F1").Value = "Sections"
$ws.Range("G1").Value = "Courses Assigned"
$ws.Range("H1").Value = "Sample Courses"
$ws.Range("I1").Value = "Utilization (Sessions/Day)"
$ws.Range("A2").Value = "C002"
$ws.Range("B2").Value = "large classroom"
$ws.Range("C2").Value = "116"
$ws.Range("D2").Value = "Projector"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = "A"
$ws.Range("G2").Value = 1
$ws.Range("H2").Value = "CS261"
$ws.Range("I2").Value = "0.4"
$ws.Range("A3").Value = "C003"
$ws.Range("B3").Value = "large classroom"
$ws.Range("C3").Value = "135"
$ws.Range("D3").Value = "Projector"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = "B"
$ws.Range("G3").Value = 1
$ws.Range("H3").Value = "CS261"
$ws.Range("I3").Value = "0.4"
$ws.Range("A4").Value = "C102"
$ws.Range("B4").Value = "classroom"
$ws.Range("C4").Value = "96"
$ws.Range("D4").Value = "Projector"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = "A, B"
$ws.Range("G4").Value = 2
$ws.Range("H4").Value = "CS263, CS264 (Tutorial)"
$ws.Range("I4").Value = "0.6"
$ws.Range("A5").Value = "C104"
$ws.Range("B5").Value = "classroom"
$ws.Range("C5").Value = "96"
$ws.Range("D5").Value = "Projector"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = "B"
$ws.Range("G5").Value = 1
$ws.Range("H5").Value = "ELECTIVE_B3"
$ws.Range("I5").Value = "0.4"
$ws.Range("A6").Value = "C201"
$ws.Range("B6").Value = "classroom"
$ws.Range("C6").Value = "96"
$ws.Range("D6").Value = "Projector"
$ws.Range("E6").Value = 6
$ws.Range("F6").Value = "A, B"
$ws.Range("G6").Value = 3
$ws.Range("H6").Value = "MA261, CS264, CS262"
$ws.Range("I6").Value = "1.2"
$ws.Range("A7").Value = "C202"
$ws.Range("B7").Value = "classroom"
$ws.Range("C7").Value = "96"
$ws.Range("D7").Value = "Projector"
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = "B"
$ws.Range("G7").Value = 1
$ws.Range("H7").Value = "MA262"
$ws.Range("I7").Value = "0.4"
$ws.Range("A8").Value = "C204"
$ws.Range("B8").Value = "classroom"
$ws.Range("C8").Value = "96"
$ws.Range("D8").Value = "Projector"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = "A, B"
$ws.Range("G8").Value = 2
$ws.Range("H8").Value = "ELECTIVE_B3, CS264 (Tutorial)"
$ws.Range("I8").Value = "0.6"
$ws.Range("A9").Value = "C303"
$ws.Range("B9").Value = "classroom"
$ws.Range("C9").Value = "96"
$ws.Range("D9").Value = "Projector"
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = "A"
$ws.Range("G9").Value = 1
$ws.Range("H9").Value = "CS263"
$ws.Range("I9").Value = "0.4"
$ws.Range("A10").Value = "C304"
$ws.Range("B10").Value = "classroom"
$ws.Range("C10").Value = "96"
$ws.Range("D10").Value = "Projector"
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = "A"
$ws.Range("G10").Value = 1
$ws.Range("H10").Value = "MA262"
$ws.Range("I10").Value = "0.4"
$ws.Range("A11").Value = "C305"
$ws.Range("B11").Value = "classroom"
$ws.Range("C11").Value = "96"
$ws.Range("D11").Value = "Projector"
$ws.Range("E11").Value = 2
$ws.Range("F11").Value = "B"
$ws.Range("G11").Value = 1
$ws.Range("H11").Value = "CS264"
$ws.Range("I11").Value = "0.4"
$ws.Range("A12").Value = "C401"
$ws.Range("B12").Value = "classroom"
$ws.Range("C12").Value = "96"
$ws.Range("D12").Value = "Projector"
$ws.Range("E12").Value = 2
$ws.Range("F12").Value = "A"
$ws.Range("G12").Value = 1
$ws.Range("H12").Value = "MA261"
$ws.Range("I12").Value = "0.4"
$ws.Range("A13").Value = "C403"
$ws.Range("B13").Value = "classroom"
$ws.Range("C13").Value = "78"
$ws.Range("D13").Value = "Projector"
$ws.Range("E13").Value = 2
$ws.Range("F13").Value = "A"
$ws.Range("G13").Value = 1
$ws.Range("H13").Value = "CS262"
$ws.Range("I13").Value = "0.4"
$ws.Range("A14").Value = "L306"
$ws.Range("B14").Value = "classroom"
$ws.Range("C14").Value = "96"
$ws.Range("D14").Value = "Computers"
$ws.Range("E14").Value = 2
$ws.Range("F14").Value = "A"
$ws.Range("G14").Value = 1
$ws.Range("H14").Value = "CS263 (Lab)"
$ws.Range("I14").Value = "0.4"
$ws.Range("A15").Value = "L406"
$ws.Range("B15").Value = "classroom"
$ws.Range("C15").Value = "78"
$ws.Range("D15").Value = "Computers"
$ws.Range("E15").Value = 2
$ws.Range("F15").Value = "B"
$ws.Range("G15").Value = 1
$ws.Range("H15").Value = "CS263 (Lab)"
$ws.Range("I15").Value = "0.4"
$ws.Range("A16").Value = "L408"
$ws.Range("B16").Value = "classroom without projector"
$ws.Range("C16").Value = "78"
$ws.Range("D16").Value = "Computers"
$ws.Range("E16").Value = 4
$ws.Range("F16").Value = "A, B"
$ws.Range("G16").Value = 1
$ws.Range("H16").Value = "CS262 (Lab)"
$ws.Range("I16").Value = "0.8"

# ===== LTPSC_Compliance =====
$ws = $wb.Worksheets.Item("LTPSC_Compliance")
$ws.Range("G2").Value = "[FAIL]"
$ws.Range("H2").Value = "[OK]"
$ws.Range("I2").Value = "[OK]"
$ws.Range("J2").Value = "[WARN] PARTIAL"
$ws.Range("G3").Value = "[OK]"
$ws.Range("H3").Value = "[OK]"
$ws.Range("I3").Value = "[OK]"
$ws.Range("J3").Value = "[OK] FULLY COMPLIANT"
$ws.Range("G4").Value = "[FAIL]"
$ws.Range("H4").Value = "[OK]"
$ws.Range("I4").Value = "[OK]"
$ws.Range("J4").Value = "[WARN] PARTIAL"
$ws.Range("G5").Value = "[FAIL]"
$ws.Range("H5").Value = "[OK]"
$ws.Range("I5").Value = "[OK]"
$ws.Range("J5").Value = "[WARN] PARTIAL"
$ws.Range("G6").Value = "[FAIL]"
$ws.Range("H6").Value = "[OK]"
$ws.Range("I6").Value = "[OK]"
$ws.Range("J6").Value = "[WARN] PARTIAL"
$ws.Range("G7").Value = "[FAIL]"
$ws.Range("H7").Value = "[FAIL]"
$ws.Range("I7").Value = "[OK]"
$ws.Range("J7").Value = "[WARN] PARTIAL"
$ws.Range("G8").Value = "[FAIL]"
$ws.Range("H8").Value = "[FAIL]"
$ws.Range("I8").Value = "[OK]"
$ws.Range("J8").Value = "[WARN] PARTIAL"

# ===== Executive_Summary =====
$ws = $wb.Worksheets.Item("Executive_Summary")
$ws.Range("C3").Value = "2025-12-12 16:58"
$ws.Range("C7").Value = "15/35"
$ws.Range("D7").Value = "Utilization: 42.9%"
$ws.Range("C9").Value = "[WARN] NEEDS REVIEW"
